$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.843.09'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '2.237.75'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.625'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.57'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +20.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.48'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.61%  '
$ws.Range('D15').Value = '2.572.62'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.08'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.65%  '
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '2.247.16'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').Value = '41.790.81'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').Value = '0.0₃0968'
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('E23').Value = '  +24.89%  '
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  +4.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.27'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('E32').Value = '  +2.34%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.21%  '
$ws.Range('E35').Value = '  +2.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +24.29%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +14.96%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.05%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0279'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.27%  '
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '69.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.38%  '
$ws.Range('E43').Value = '  +16.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.64'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +22.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +12.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('E48').Value = '  +3.05%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.43%  '
$ws.Range('B51').Value = 'BitTorrent-New'
$ws.Range('C51').Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range('D51').Value = '0.0₃0155'
$ws.Range('E51').Value = '  +22.08%  '
